# "Generate Report for Handoff"
# Updates the localization-status report: the zh-cn / de-de items have moved
# from "In Translation" to "Ready for handoff", the corresponding handoff
# timestamps are refreshed, and the status columns are re-widened so the new
# (longer) status text isn't truncated.

$wb = $excel.ActiveWorkbook

# ----- Overview sheet ------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = "Ready for handoff"
$ov.Range("F2").Value = "Ready for handoff"
$ov.Range("G2").Value = "2016-08-19 17:02:19"
$ov.Columns.Item(5).ColumnWidth = 16.33
$ov.Columns.Item(6).ColumnWidth = 16.33

# ----- zh-cn sheet -----------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("H2").Value = "2016-08-19 17:02:15"
$zh.Columns.Item(3).ColumnWidth = 16.33

# ----- de-de sheet -----------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = "Ready for handoff"
$de.Range("H2").Value = "2016-08-19 17:02:19"
$de.Columns.Item(3).ColumnWidth = 16.33
